$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 04.02.2022 13:15"

# Tesco row (row 3): prices swapped, delta now a negative-string, date now a text timestamp
$ws.Range("B3").Value = 35.5
$ws.Range("C3").Value = 35.51
$ws.Range("D3").Value = "-0.01"
$ws.Range("E3").Style = "Normal"
$ws.Range("E3").Value = "2022-02-04 13:15:11"
